# Adds a new daily log sheet "2024-05-31" (after "2024-05-30") and
# appends the corresponding summary row to the "current" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new daily sheet, placed right after the last existing
#        daily sheet ("2024-05-30"), and name it "2024-05-31". ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2024-05-31"

# Header row, matching the layout used by the other daily sheets.
$newSheet.Range("A1").Value = "Время"
$newSheet.Range("B1").Value = "ФИО пациента"
$newSheet.Range("C1").Value = "М\Ж\Р"
$newSheet.Range("D1").Value = "Дата рождения"
$newSheet.Range("E1").Value = "Причина"
$newSheet.Range("F1").Value = "Давление"

# --- 2. Append the summary row for 2024-05-31 to the "current" sheet. -
$currentSheet = $wb.Worksheets.Item("current")

# A11 / D11 hold digit-like text ("2024-05-31", "3") that must stay text
# rather than being auto-converted to a date/number, so the cell is
# briefly formatted as Text, written, then the format override is
# cleared again to match the rest of the sheet (which carries no
# explicit per-cell style).
$cellA11 = $currentSheet.Range("A11")
$cellA11.NumberFormat = "@"
$cellA11.Value = "2024-05-31"
$cellA11.ClearFormats()

$currentSheet.Range("B11").Value = "Karp_Kuzmin"
$currentSheet.Range("C11").Value = -1

$cellD11 = $currentSheet.Range("D11")
$cellD11.NumberFormat = "@"
$cellD11.Value = "3"
$cellD11.ClearFormats()

# --- 3. Restore the original active sheet/selection ("2024-05-30"!B13)
#        so that adding the new sheet doesn't steal the tab focus. ----
$lastSheet.Activate() | Out-Null
$lastSheet.Range("B13").Select() | Out-Null
